# Add columns I (I0) and J (IF) with header + data values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header-cell formatting (bold font, border, centered alignment) from H1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data values: row -> (I value, J value)
$data = @{
    2  = @(4, 4)
    3  = @(3, 3)
    4  = @(7, 8)
    5  = @(9, 9)
    6  = @(5, 6)
    7  = @(9, 9)
    8  = @(5, 6)
    9  = @(6, 6)
    10 = @(8, 8)
    11 = @(8, 9)
    12 = @(8, 9)
    13 = @(6, 7)
    14 = @(5, 6)
    15 = @(8, 8)
    16 = @(5, 6)
    17 = @(4, 4)
    18 = @(7, 7)
    19 = @(6, 6)
    20 = @(4, 5)
    21 = @(8, 8)
    22 = @(7, 8)
    23 = @(7, 8)
    24 = @(8, 8)
    25 = @(6, 7)
    26 = @(12, 13)
    27 = @(2, 2)
    28 = @(6, 6)
    29 = @(8, 8)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
